# Reorder the "Recorded By" contributor lists in column G of the
# "Session Analysis Results" sheet so that "System"/"system" entries
# come first, followed by the remaining contributor(s) (sorted
# alphabetically when there is more than one non-System contributor).
#
# This mirrors the upstream sync that re-serialized the recorded-by
# list with System listed first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Value = "System, system, backup@backdoor.com"; Rows = @(2, 28, 54)},
    @{Value = "System, dnasr281@gmail.com"; Rows = @(3, 6, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26, 29, 32, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52, 55, 58, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 93, 94, 96, 99, 101, 109, 110, 111, 112, 116, 118, 119, 120, 122, 125, 127, 135, 136, 137, 138, 142, 144, 145, 146, 148, 151, 153)},
    @{Value = "System, backup@backdoor.com"; Rows = @(4, 5, 8, 30, 31, 34, 56, 57, 60, 80, 81, 82, 106, 107, 108, 132, 133, 134)},
    @{Value = "System, admin@admin.com"; Rows = @(7, 33, 59)},
    @{Value = "admin@admin.com, dnasr281@gmail.com"; Rows = @(87, 113, 139)}
)

foreach ($entry in $updates) {
    foreach ($row in $entry.Rows) {
        $ws.Range("G$row").Value = $entry.Value
    }
}

Write-Host "Updated $(($updates | ForEach-Object { $_.Rows.Count } | Measure-Object -Sum).Sum) cells in column G"
